$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "Performance fee paid" header column (column N). Excel shifts
# the following column ("GPS ITD") left into N, updates the used
# dimension/row span, and drops the column from the data-validation sqref
# list automatically.
$ws.Columns.Item(14).Delete()

# Reset the view: clear the scrolled-right position and select A3.
$ws.Range("A3").Select() | Out-Null

# Lock the workbook structure (sheets can't be added/moved/deleted/hidden).
$wb.Protect($null, $true, $false)

# Re-protect the worksheet so the previously stored password hash is
# cleared while leaving the sheet protected.
$ws.Protect()
